{"js": "// The template paragraph (centered, underlined, Amasis MT Pro Black, 18pt)\n// already holding {{NAME}} is the last paragraph in the body. Insert four\n// new paragraphs - {{NAME}}, {{SOCIAL}}, {{EMPTY}}, {{EMPTY}} - right before\n// it. Paragraph.insertParagraph(text, \"Before\") clones the anchor\n// paragraph's paragraph/run formatting (pPr/rPr), so each new paragraph\n// automatically picks up the same centered / underlined / font / size\n// formatting as the existing {{NAME}} paragraph.\nconst texts = [\"{{NAME}}\", \"{{SOCIAL}}\", \"{{EMPTY}}\", \"{{EMPTY}}\"];\n\nfor (const t of texts) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const anchor = paragraphs.items[paragraphs.items.length - 1];\n  anchor.insertParagraph(t, \"Before\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The template paragraph (centered, underlined, Amasis MT Pro Black, 18pt)\n# already holding {{NAME}} is the last paragraph in the body. We insert four\n# new paragraphs - {{NAME}}, {{SOCIAL}}, {{EMPTY}}, {{EMPTY}} - right before\n# it, each inheriting that same paragraph/run formatting (InsertParagraphBefore\n# clones the pPr/rPr of the anchor paragraph), then fill in the text.\n$texts = @(\"{{NAME}}\", \"{{SOCIAL}}\", \"{{EMPTY}}\", \"{{EMPTY}}\")\n\nforeach ($t in $texts) {\n    $lastIndex = $d.Paragraphs.Count\n    $anchor = $d.Paragraphs.Item($lastIndex)\n    $anchor.Range.InsertParagraphBefore()\n    $newPara = $d.Paragraphs.Item($lastIndex)\n    $newPara.Range.Text = $t\n}\n"}
